# Weekly update: insert a new "Orégano" price observation as row 138,
# pushing the existing rows 138-146 down to 139-147 (dimension grows to R147).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138; Excel shifts rows 138:146 down to 139:147.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new weekly record.
$ws.Cells.Item(138, 1).Value  = 6
$ws.Cells.Item(138, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(138, 3).Value  = "Metropolitana"
$ws.Cells.Item(138, 4).Value  = 44585
$ws.Cells.Item(138, 5).Value  = 13
$ws.Cells.Item(138, 6).Value  = 100112029
$ws.Cells.Item(138, 7).Value  = "Orégano"
$ws.Cells.Item(138, 8).Value  = "Sin especificar"
$ws.Cells.Item(138, 9).Value  = "Primera"
$ws.Cells.Item(138, 10).Value = 29
$ws.Cells.Item(138, 11).Value = 8000
$ws.Cells.Item(138, 12).Value = 9000
$ws.Cells.Item(138, 13).Value = 8483
$ws.Cells.Item(138, 14).Value = "`$/docena de atados"
$ws.Cells.Item(138, 15).Value = "Región Metropolitana"
$ws.Cells.Item(138, 16).Value = 2828
$ws.Cells.Item(138, 17).Value = 3
$ws.Cells.Item(138, 18).Value = "Hortaliza"

# Apply the same date number format (style used by column D) to the new cell.
$ws.Cells.Item(138, 4).NumberFormat = $ws.Cells.Item(139, 4).NumberFormat
